# options_guide.xlsx — "Fix reduce COLA/Chain-CPI COLA"
#
# The BPC file-shell naming scheme was reworked (several *.xlsx shell
# filenames were renamed/shortened) and the mismatched option labels for
# the "COLAChainCPI" / "6-reduceCOLA" rows were corrected/swapped.
# D2:D39 are CONCATENATE(E,F) formulas, so they recompute automatically
# once E/F are updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Opt0 / equivalent ------------------------------------------------
$ws.Range("F3").Value2 = "BPCtableShellsEQOPT0.xlsx"

# --- Row 4: Opt0 / Payable Law / per capita ---------------------------------
$ws.Range("F4").Value2 = "BPCtableShellsPayable.xlsx"

# --- Row 5: Opt0 / Payable Law / equivalent ---------------------------------
$ws.Range("F5").Value2 = "BPCtableShellsEquivalentPayable.xlsx"

# --- Row 6: 1-mini.pia / Annual PIA / per capita ----------------------------
$ws.Range("F6").Value2 = "BPCtableShellsMiniPIA.xlsx"

# --- Row 7: 1-mini.pia / Annual PIA / equivalent ----------------------------
$ws.Range("F7").Value2 = "BPCtableShellsEquivalentMiniPIA.xlsx"

# --- Row 10: 3-capSpouse / Cap Spouse Benefits / per capita -----------------
$ws.Range("F10").Value2 = "BPCtableShellsCapSpouse.xlsx"

# --- Rows 20-21: COLAChainCPI -- mislabeled "Full Chained-CPI COLA" --------
#     fixed to "Reduce COLA"
$ws.Range("B20").Value2 = "Reduce COLA"
$ws.Range("B21").Value2 = "Reduce COLA"

# --- Rows 22-23: 6-reduceCOLA -- mislabeled "Partial Chained-CPI COLA" -----
#     fixed to "Chained-CPI COLA"
$ws.Range("B22").Value2 = "Chained-CPI COLA"
$ws.Range("B23").Value2 = "Chained-CPI COLA"

# --- Rows 28-29: 12-taxmax150000 directory now includes the option suffix --
$ws.Range("E28").Value2 = "X:\programs\run912\opt12(taxmax150000)\"
$ws.Range("E29").Value2 = "X:\programs\run912\opt12(taxmax150000)\"

# --- Window state: selection moved to B24 (view scrolled so row 4 is top) --
$ws.Range("B24").Select()

$wb.Save()
